$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineShape($headerFooter, $newName) {
    $shp = $headerFooter.Range.InlineShapes.Item(1)
    # Routing the rename through the Selection (rather than the raw
    # InlineShapes collection handle) avoids a stale-handle snag that the
    # header/footer story can hit when addressed directly.
    $shp.Range.Select()
    $sel = $word.Selection
    $sel.InlineShapes.Item(1).Name = $newName
}

# Pearson logo (footers) : image2.png -> image1.png
Rename-InlineShape $sec.Footers(1) "image1.png"
Rename-InlineShape $sec.Footers(2) "image1.png"

# BTEC logo (headers) : image1.jpg -> image2.jpg
Rename-InlineShape $sec.Headers(1) "image2.jpg"
Rename-InlineShape $sec.Headers(2) "image2.jpg"

Write-Host "renamed all 4 inline pictures"
